# Update the weekly fruit/vegetable price data: swap the row 3 and row 4
# values for Date, Volumen, Precio minimo, Precio maximo, Precio promedio
# ponderado and Precio $/Kg columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (becomes the values previously held by row 4)
$ws.Range("D3").Value = 44414
$ws.Range("J3").Value = 500
$ws.Range("K3").Value = 31000
$ws.Range("L3").Value = 32000
$ws.Range("M3").Value = 31500
$ws.Range("P3").Value = 1260

# Row 4 (becomes the values previously held by row 3)
$ws.Range("D4").Value = 44827
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 30000
$ws.Range("L4").Value = 31000
$ws.Range("M4").Value = 30500
$ws.Range("P4").Value = 1220
